# Applies the benchmark-result update described in the commit diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- Simple single-value cell updates (unique text, safe to Find/Replace) ---
Replace-Text "89.91" "0M"
Replace-Text "191.13" "0M"
Replace-Text "1895" "0M"
Replace-Text "25834" "25837"
Replace-Text "0.04270" "0.04282"
Replace-Text "0.00523" "0.00517"
Replace-Text "190.96827" "191.13064"

# --- Collapse the three multi-run/tab-separated cells to single values ---
$t = $d.Tables.Item(1)
$t.Cell(44, 1).Range.Text = "89.91"
$t.Cell(45, 1).Range.Text = "191.13"
$t.Cell(46, 1).Range.Text = "1895"
